$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" footer timestamp shown in cell A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 22:22"

# Update country rows: names (col A) and statistics (cols B:H) per the new data snapshot.
# Some countries were re-sorted (moved to a new row position) and/or had updated counts.
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Range("B4").Value = 422510
$ws.Range("C4").Value = 22175
$ws.Range("D4").Value = 22187
$ws.Range("E4").Value = 385850
$ws.Range("F4").Value = 9233
$ws.Range("G4").Value = 1632
$ws.Range("H4").Value = 14473

$ws.Cells.Item(16, 1).Value = "Canada"
$ws.Range("B16").Value = 19195
$ws.Range("C16").Value = 1298
$ws.Range("D16").Value = 4533
$ws.Range("E16").Value = 14235
$ws.Range("F16").Value = 426
$ws.Range("G16").Value = 46
$ws.Range("H16").Value = 427

$ws.Cells.Item(17, 1).Value = "Brasil"
$ws.Range("B17").Value = 15927
$ws.Range("C17").Value = 1893
$ws.Range("D17").Value = 127
$ws.Range("E17").Value = 15000
$ws.Range("F17").Value = 296
$ws.Range("G17").Value = 114
$ws.Range("H17").Value = 800

$ws.Cells.Item(19, 1).Value = "Austria"
$ws.Range("B19").Value = 12941
$ws.Range("C19").Value = 302
$ws.Range("D19").Value = 4512
$ws.Range("E19").Value = 8156
$ws.Range("F19").Value = 267
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 273

$ws.Cells.Item(75, 1).Value = "Camerun"
$ws.Range("B75").Value = 730
$ws.Range("C75").Value = 45
$ws.Range("D75").Value = 60
$ws.Range("E75").Value = 660
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 10

$ws.Cells.Item(76, 1).Value = "Kazajistan"
$ws.Range("B76").Value = 727
$ws.Range("C76").Value = 30
$ws.Range("D76").Value = 54
$ws.Range("E76").Value = 666
$ws.Range("F76").Value = 21
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 7

$ws.Cells.Item(77, 1).Value = "Crucero"
$ws.Range("B77").Value = 712
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 619
$ws.Range("E77").Value = 82
$ws.Range("F77").Value = 10
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 11

$ws.Cells.Item(79, 1).Value = "Tunez"
$ws.Range("B79").Value = 628
$ws.Range("C79").Value = 5
$ws.Range("D79").Value = 25
$ws.Range("E79").Value = 579
$ws.Range("F79").Value = 67
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 24

$ws.Cells.Item(94, 1).Value = "Costa de Marfil"
$ws.Range("B94").Value = 384
$ws.Range("C94").Value = 35
$ws.Range("D94").Value = 48
$ws.Range("E94").Value = 333
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 3

$ws.Cells.Item(95, 1).Value = "Taiwan"
$ws.Range("B95").Value = 379
$ws.Range("C95").Value = 3
$ws.Range("D95").Value = 67
$ws.Range("E95").Value = 307
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 5

$ws.Cells.Item(96, 1).Value = "Reunion"
$ws.Range("B96").Value = 358
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 40
$ws.Range("E96").Value = 318
$ws.Range("F96").Value = 4
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0

$ws.Cells.Item(97, 1).Value = "Jordania"
$ws.Range("B97").Value = 358
$ws.Range("C97").Value = 5
$ws.Range("D97").Value = 150
$ws.Range("E97").Value = 202
$ws.Range("F97").Value = 5
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 6

$ws.Cells.Item(103, 1).Value = "Nigeria"
$ws.Range("B103").Value = 276
$ws.Range("C103").Value = 22
$ws.Range("D103").Value = 44
$ws.Range("E103").Value = 226
$ws.Range("F103").Value = 2
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 6

$ws.Cells.Item(104, 1).Value = "Mauricio"
$ws.Range("B104").Value = 273
$ws.Range("C104").Value = 5
$ws.Range("D104").Value = 19
$ws.Range("E104").Value = 247
$ws.Range("F104").Value = 3
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 7

$ws.Cells.Item(105, 1).Value = "Kirguistan"
$ws.Range("B105").Value = 270
$ws.Range("C105").Value = 42
$ws.Range("D105").Value = 33
$ws.Range("E105").Value = 233
$ws.Range("F105").Value = 5
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 4

$ws.Cells.Item(106, 1).Value = "Estado de Palestina"
$ws.Range("B106").Value = 263
$ws.Range("C106").Value = 2
$ws.Range("D106").Value = 44
$ws.Range("E106").Value = 218
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 1

$ws.Cells.Item(120, 1).Value = "Martinica"
$ws.Range("B120").Value = 154
$ws.Range("C120").Value = 2
$ws.Range("D120").Value = 50
$ws.Range("E120").Value = 98
$ws.Range("F120").Value = 19
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 6

$ws.Cells.Item(122, 1).Value = "Guadalupe"
$ws.Range("B122").Value = 141
$ws.Range("C122").Value = 2
$ws.Range("D122").Value = 43
$ws.Range("E122").Value = 90
$ws.Range("F122").Value = 13
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 8

$ws.Cells.Item(144, 1).Value = "Congo"
$ws.Range("B144").Value = 45
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 2
$ws.Range("E144").Value = 38
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 5

$ws.Cells.Item(145, 1).Value = "Islas Caimanes"
$ws.Range("B145").Value = 45
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 6
$ws.Range("E145").Value = 38
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 1

$ws.Cells.Item(153, 1).Value = "Eritrea"
$ws.Range("B153").Value = 33
$ws.Range("C153").Value = 2
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 33
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 0

$ws.Cells.Item(154, 1).Value = "Guyana"
$ws.Range("B154").Value = 33
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 8
$ws.Range("E154").Value = 19
$ws.Range("F154").Value = 4
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 6

$ws.Cells.Item(155, 1).Value = "Guam"
$ws.Range("B155").Value = 32
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 0
$ws.Range("E155").Value = 31
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 1

$ws.Cells.Item(156, 1).Value = "San Martin (Parte Francesa)"
$ws.Range("B156").Value = 32
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 7
$ws.Range("E156").Value = 23
$ws.Range("F156").Value = 6
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 2

$ws.Cells.Item(165, 1).Value = "Angola"
$ws.Range("B165").Value = 19
$ws.Range("C165").Value = 2
$ws.Range("D165").Value = 2
$ws.Range("E165").Value = 15
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 2

$ws.Cells.Item(166, 1).Value = "Siria"
$ws.Range("B166").Value = 19
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 4
$ws.Range("E166").Value = 13
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 2

$ws.Cells.Item(167, 1).Value = "Maldivas"
$ws.Range("B167").Value = 19
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 13
$ws.Range("E167").Value = 6
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

$ws.Cells.Item(168, 1).Value = "Nueva Caledonia"
$ws.Range("B168").Value = 18
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 1
$ws.Range("E168").Value = 17
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

$ws.Cells.Item(169, 1).Value = "Guinea Ecuatorial"
$ws.Range("B169").Value = 18
$ws.Range("C169").Value = 2
$ws.Range("D169").Value = 3
$ws.Range("E169").Value = 15
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0

$ws.Cells.Item(170, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B170").Value = 17
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 0
$ws.Range("E170").Value = 17
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

$ws.Cells.Item(171, 1).Value = "Mozambique"
$ws.Range("B171").Value = 17
$ws.Range("C171").Value = 7
$ws.Range("D171").Value = 1
$ws.Range("E171").Value = 16
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

$ws.Cells.Item(174, 1).Value = "Fiyi"
$ws.Range("B174").Value = 15
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 15
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

$ws.Cells.Item(175, 1).Value = "Laos"
$ws.Range("B175").Value = 15
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 15
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

$ws.Cells.Item(179, 1).Value = "Curazao"
$ws.Range("B179").Value = 14
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 7
$ws.Range("E179").Value = 6
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 1
